# Mozilla Firefox.pptx -- "powerpoint changes -- benefits"
#
# The original slide 9 ("Critical Analysis") is duplicated. The duplicate is
# moved to sit BEFORE the original (so it keeps slide position 9) and becomes
# the "Critical Analysis -- Drawbacks" slide (essentially the old content,
# lightly reworded). The original slide (now at position 10) is rewritten
# into the new "Critical Analysis -- Benefits" slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Duplicate the "Critical Analysis" slide (slide 9) and move the new
#    copy in front of the original, so that:
#      position 9  -> new copy       -> becomes "Drawbacks"
#      position 10 -> original slide -> becomes "Benefits"
# ---------------------------------------------------------------------
$original = $p.Slides.Item(9)
$dup = $original.Duplicate()
$dup.MoveTo(9)

# ---------------------------------------------------------------------
# 2. Rewrite the duplicate (now at position 9) into the "Drawbacks" slide.
# ---------------------------------------------------------------------
$drawbacks = $p.Slides.Item(9)

$dTitle = $drawbacks.Shapes.Item(1).TextFrame.TextRange
$dTitle.Text = "Critical Analysis -- Drawbacks"

$dBody = $drawbacks.Shapes.Item(2).TextFrame.TextRange
$dBody.Text = "In some areas, Firefox can be difficult to maintain`r" + `
    "A lot of dependencies throughout the browser`r" + `
    "If one API or library breaks, the browser could break`r" + `
    "`r" + `
    "Openness can lead to performance issues`r" + `
    "Extensions are easily developed by users and uploaded to the Firefox Store`r" + `
    "Open Nature can lead to user-created extensions bogging down browser performance`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    "`r" + `
    ""

for ($i = 2; $i -le 4; $i++) { $dBody.Paragraphs($i, 1).IndentLevel = 2 }
for ($i = 6; $i -le 16; $i++) { $dBody.Paragraphs($i, 1).IndentLevel = 2 }

# ---------------------------------------------------------------------
# 3. Rewrite the original slide (now at position 10) into the "Benefits"
#    slide.
# ---------------------------------------------------------------------
$benefits = $p.Slides.Item(10)

$bTitle = $benefits.Shapes.Item(1).TextFrame.TextRange
$bTitle.Text = "Critical Analysis -- Benefits"

$bShape = $benefits.Shapes.Item(2)
$bBody = $bShape.TextFrame.TextRange
$dash = [char]0x2013
$apos = [char]0x2019
$bBody.Text = "Firefox is Platform Independent `r" + `
    "Due to its architecture" + $apos + "s reusable components, Firefox has been ported to and is available on most popular Operating Systems, leading to an often consistent user experience.`r" + `
    "Firefox is maintainable`r" + `
    "Layered Architecture`r" + `
    "Components are restricted to only communicate with certain layers. Changes don" + $apos + "t usually affect the entire code base.`r" + `
    "Firefox is extensible`r" + `
    "Many extensions are available`r" + `
    "Extensions are independently built " + $dash + " can be easily installed and removed`r" + `
    "Extensions can add a lot of interesting functionality not initially available`r" + `
    ""

$bBody.Paragraphs(2, 1).IndentLevel = 2
for ($i = 4; $i -le 5; $i++) { $bBody.Paragraphs($i, 1).IndentLevel = 2 }
for ($i = 7; $i -le 9; $i++) { $bBody.Paragraphs($i, 1).IndentLevel = 2 }

# Apply the Wingdings symbol font to the leading "C" of the "Components..."
# bullet, matching the original author's formatting.
$cRun = $bBody.Paragraphs(5, 1).Characters(1, 1)
$cRun.Font.Name = "Wingdings"

# Shrink text slightly so the (now much longer) body still fits the
# placeholder -- mirrors PowerPoint's automatic "Shrink text on overflow".
$bShape.TextFrame2.AutoSize = 2
$bShape.TextFrame2.FontScale = 90
$bShape.TextFrame2.LineSpaceReduction = 10
